# Modificación de fechas en presentaciones.
# Update the cover-slide subtitle date from "Mayo 2012" to "Abril 2013".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

foreach ($shape in $s.Shapes) {
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -eq "Mayo 2012") {
            $tr.Text = "Abril 2013"
        }
    }
}
